$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "tsv initial format fix": the shared string that used to read "PSEUDO_HEARING"
# should just read "HEARING" (matching the already-existing HEARING value used
# elsewhere in the sheet, e.g. C7).
$ws.Range("C6").Value = "HEARING"

# Author also left the active selection on C6 (previously on C11), and moved the
# tab-ratio (sheet-tab-area/horizontal-scrollbar split) slider - apply both.
[void]$ws.Range("C6").Select()

$win = $wb.Windows.Item(1)
$win.TabRatio = 0.5
